$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.099.96"
$ws.Range("E2").Value = "  +2.60%  "
$ws.Range("D3").Value = "2.403.25"
$ws.Range("E3").Value = "  +2.99%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'560.42"
$ws.Range("E5").Value = "  +2.44%  "
$ws.Range("D6").Value = "'138.14"
$ws.Range("E6").Value = "  +4.94%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  +0.86%  "
$ws.Range("D9").Value = "2.402.04"
$ws.Range("E9").Value = "  +3.09%  "
$ws.Range("E10").Value = "  +2.81%  "
$ws.Range("D11").Value = "'5.70"
$ws.Range("E11").Value = "  +3.56%  "
$ws.Range("E12").Value = "  +0.11%  "
$ws.Range("E13").Value = "  +3.59%  "
$ws.Range("D14").Value = "'25.68"
$ws.Range("E14").Value = "  +7.72%  "
$ws.Range("D15").Value = "2.831.66"
$ws.Range("E15").Value = "  +2.98%  "
$ws.Range("D16").Value = "62.053.25"
$ws.Range("E16").Value = "  +2.58%  "
$ws.Range("E17").Value = "  +3.69%  "
$ws.Range("D18").Value = "2.413.36"
$ws.Range("E18").Value = "  +3.42%  "
$ws.Range("D19").Value = "'11.00"
$ws.Range("E19").Value = "  +3.40%  "
$ws.Range("D20").Value = "'343.70"
$ws.Range("E20").Value = "  +9.02%  "
$ws.Range("D21").Value = "'4.22"
$ws.Range("E21").Value = "  +1.41%  "
$ws.Range("E22").Value = "  +3.51%  "
$ws.Range("E23").Value = "  +0.39%  "
$ws.Range("D24").Value = "'65.00"
$ws.Range("E24").Value = "  +1.27%  "
$ws.Range("E25").Value = "  +1.74%  "
$ws.Range("E26").Value = "  +0.23%  "
$ws.Range("D27").Value = "'8.33"
$ws.Range("E27").Value = "  +5.45%  "
$ws.Range("D28").Value = "'1.50"
$ws.Range("E28").Value = "  +10.89%  "
$ws.Range("D29").Value = "'1.37"
$ws.Range("E29").Value = "  +14.12%  "
$ws.Range("E30").Value = "  +4.00%  "
$ws.Range("D31").Value = "0.0₃0771"
$ws.Range("E31").Value = "  +4.71%  "
$ws.Range("D32").Value = "'6.37"
$ws.Range("D33").Value = "'171.60"
$ws.Range("E33").Value = "  -0.86%  "
$ws.Range("D34").Value = "'1.40"
$ws.Range("E34").Value = "  +1.59%  "
$ws.Range("E35").Value = "  +3.35%  "
$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D36").Value = "'4.53"
$ws.Range("E36").Value = "  +10.84%  "
$ws.Range("B37").Value = "EthereumClassic"
$ws.Range("C37").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D37").Value = "'18.50"
$ws.Range("E37").Value = "  +3.37%  "
$ws.Range("B38").Value = "USDe"
$ws.Range("C38").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D38").Value = "'0.998"
$ws.Range("E38").Value = "  -0.05%  "
$ws.Range("B39").Value = "Bittensor"
$ws.Range("C39").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D39").Value = "'357.08"
$ws.Range("E39").Value = "  +8.61%  "
$ws.Range("E40").Value = "  -0.07%  "
$ws.Range("D41").Value = "'1.68"
$ws.Range("E41").Value = "  +8.63%  "
$ws.Range("D42").Value = "'39.03"
$ws.Range("E42").Value = "  +2.76%  "
$ws.Range("D43").Value = "'143.34"
$ws.Range("E43").Value = "  +2.90%  "
$ws.Range("D44").Value = "'3.66"
$ws.Range("E44").Value = "  +4.76%  "
$ws.Range("D45").Value = "'20.36"
$ws.Range("E45").Value = "  +5.13%  "
$ws.Range("D46").Value = "'0.0964"
$ws.Range("E46").Value = "  +1.97%  "
$ws.Range("E47").Value = "  +4.14%  "
$ws.Range("E48").Value = "  +3.50%  "
$ws.Range("E49").Value = "  +3.46%  "
$ws.Range("D50").Value = "'17.79"
$ws.Range("E50").Value = "  +4.52%  "
$ws.Range("D51").Value = "0.0₆0218"
$ws.Range("E51").Value = "  -0.56%  "
